# Remove column A (the stray/misaligned leading column) and shift the
# remaining columns B:F left into A:E, so that the header row and data
# rows line up correctly under columns A-E.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A:A").Delete()
